# Auto-generated Excel COM-interop script
# Applies scheduled market-data refresh values to the profit sheets
# (columns H-N: currentAveragePrice*, LevePrice*, LeveProfit*)

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1701.4222
$ws.Range("J17").Value = 1764.1
$ws.Range("L17").Value = 5292.299999999999
$ws.Range("N17").Value = -5628.299999999999
$ws.Range("H32").Value = 1229.6364
$ws.Range("I32").Value = 871.5
$ws.Range("J32").Value = 1659.4
$ws.Range("K32").Value = 871.5
$ws.Range("L32").Value = 1659.4
$ws.Range("M32").Value = -545.5
$ws.Range("N32").Value = -2311.4
$ws.Range("H44").Value = 2000
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()
$ws.Range("H64").Value = 7454.9375
$ws.Range("I64").Value = 4780.4
$ws.Range("J64").Value = 8670.637000000001
$ws.Range("K64").Value = 4780.4
$ws.Range("L64").Value = 8670.637000000001
$ws.Range("M64").Value = -4532.4
$ws.Range("N64").Value = -9166.637000000001
$ws.Range("H67").Value = 7454.9375
$ws.Range("I67").Value = 4780.4
$ws.Range("J67").Value = 8670.637000000001
$ws.Range("K67").Value = 4780.4
$ws.Range("L67").Value = 8670.637000000001
$ws.Range("M67").Value = -3922.4
$ws.Range("N67").Value = -10386.637
$ws.Range("H113").Value = 3449
$ws.Range("J113").Value = 4699
$ws.Range("L113").Value = 4699
$ws.Range("N113").Value = -11207
$ws.Range("H116").Value = 56304.57
$ws.Range("I116").Value = 82112.78
$ws.Range("K116").Value = 82112.78
$ws.Range("M116").Value = -78670.78
$ws.Range("H129").Value = 1691.6666
$ws.Range("I129").Value = 1443.3334
$ws.Range("J129").Value = 2933.3333
$ws.Range("K129").Value = 4330.0002
$ws.Range("L129").Value = 8799.999899999999
$ws.Range("M129").Value = 669.9997999999996
$ws.Range("N129").Value = -18799.9999
$ws.Range("H132").Value = 53495.977
$ws.Range("I132").Value = 55891.14
$ws.Range("K132").Value = 167673.42
$ws.Range("M132").Value = -165143.42

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5620248.5
$ws.Range("I32").Value = 5620248.5
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 5620248.5
$ws.Range("L32").Value = 0
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -5619961.5
$ws.Range("H88").Value = 1793.625
$ws.Range("I88").Value = 1570
$ws.Range("K88").Value = 1570
$ws.Range("M88").Value = -1164
$ws.Range("H91").Value = 1793.625
$ws.Range("I91").Value = 1570
$ws.Range("K91").Value = 1570
$ws.Range("M91").Value = -166
$ws.Range("H122").Value = 4999.5
$ws.Range("I122").Value = 4998.6665
$ws.Range("K122").Value = 14995.9995
$ws.Range("M122").Value = -12545.9995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 2569.375
$ws.Range("I64").Value = 305
$ws.Range("K64").Value = 305
$ws.Range("M64").Value = -80
$ws.Range("H67").Value = 2569.375
$ws.Range("I67").Value = 305
$ws.Range("K67").Value = 305
$ws.Range("M67").Value = 475
$ws.Range("H86").Value = 2268.5334
$ws.Range("I86").Value = 2411
$ws.Range("J86").Value = 2054.8333
$ws.Range("K86").Value = 2411
$ws.Range("L86").Value = 2054.8333
$ws.Range("M86").Value = -1288
$ws.Range("N86").Value = -4300.8333
$ws.Range("H89").Value = 2268.5334
$ws.Range("I89").Value = 2411
$ws.Range("J89").Value = 2054.8333
$ws.Range("K89").Value = 12055
$ws.Range("L89").Value = 10274.1665
$ws.Range("M89").Value = -6439
$ws.Range("N89").Value = -21506.1665
$ws.Range("H94").Value = 1608.1538
$ws.Range("I94").Value = 1064.1578
$ws.Range("K94").Value = 1064.1578
$ws.Range("M94").Value = -613.1578
$ws.Range("H107").Value = 5082.0835
$ws.Range("I107").Value = 4788.778
$ws.Range("J107").Value = 5962
$ws.Range("K107").Value = 4788.778
$ws.Range("L107").Value = 5962
$ws.Range("M107").Value = -2868.778
$ws.Range("N107").Value = -9802
$ws.Range("H134").Value = 499604.5
$ws.Range("I134").Value = 606666.6
$ws.Range("J134").Value = 11877.111
$ws.Range("K134").Value = 1819999.8
$ws.Range("L134").Value = 35631.333
$ws.Range("M134").Value = -1817464.8
$ws.Range("N134").Value = -40701.333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("N45").ClearContents()
$ws.Range("H88").Value = 14068.4
$ws.Range("J88").Value = 14068.4
$ws.Range("L88").Value = 14068.4
$ws.Range("N88").Value = -14880.4
$ws.Range("H91").Value = 14068.4
$ws.Range("J91").Value = 14068.4
$ws.Range("L91").Value = 14068.4
$ws.Range("N91").Value = -16876.4
$ws.Range("H122").Value = 4270.3335
$ws.Range("J122").Value = 5389.125
$ws.Range("L122").Value = 16167.375
$ws.Range("N122").Value = -21067.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H118").Value = 1076.3334
$ws.Range("I118").Value = 1064.5
$ws.Range("J118").Value = 1100
$ws.Range("K118").Value = 3193.5
$ws.Range("L118").Value = 3300
$ws.Range("M118").Value = -1950.5
$ws.Range("N118").Value = -5786
$ws.Range("H122").Value = 695368.6
$ws.Range("J122").Value = 1069529.8
$ws.Range("L122").Value = 9625768.200000001
$ws.Range("N122").Value = -9630668.200000001
$ws.Range("H139").Value = 4418
$ws.Range("I139").Value = 5000
$ws.Range("J139").Value = 4127
$ws.Range("K139").Value = 15000
$ws.Range("L139").Value = 12381
$ws.Range("M139").Value = -9860
$ws.Range("N139").Value = -22661

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 4665
$ws.Range("J12").Value = 4500
$ws.Range("L12").Value = 4500
$ws.Range("N12").Value = -4780
$ws.Range("H44").Value = 50000
$ws.Range("J44").Value = 50000
$ws.Range("L44").Value = 50000
$ws.Range("N44").Value = -51192
$ws.Range("H102").Value = 1613.6904
$ws.Range("I102").Value = 1186.2354
$ws.Range("J102").Value = 3430.375
$ws.Range("K102").Value = 1186.2354
$ws.Range("L102").Value = 3430.375
$ws.Range("M102").Value = 435.7646
$ws.Range("N102").Value = -6674.375
$ws.Range("H104").Value = 100671
$ws.Range("J104").Value = 100671
$ws.Range("L104").Value = 100671
$ws.Range("N104").Value = -107659

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2923
$ws.Range("I7").Value = 2718.5715
$ws.Range("K7").Value = 2718.5715
$ws.Range("M7").Value = -2606.5715
$ws.Range("H9").Value = 1175.9
$ws.Range("I9").Value = 420
$ws.Range("J9").Value = 1931.8
$ws.Range("K9").Value = 420
$ws.Range("L9").Value = 1931.8
$ws.Range("M9").Value = -196
$ws.Range("N9").Value = -2379.8
$ws.Range("H16").Value = 575.5714
$ws.Range("I16").Value = 583.73334
$ws.Range("J16").Value = 555.1667
$ws.Range("K16").Value = 583.73334
$ws.Range("L16").Value = 555.1667
$ws.Range("M16").Value = -413.73334
$ws.Range("N16").Value = -895.1667
$ws.Range("H40").Value = 5567.25
$ws.Range("I40").Value = 4648.357
$ws.Range("K40").Value = 4648.357
$ws.Range("M40").Value = -4512.357
$ws.Range("H55").Value = 1129.1818
$ws.Range("I55").Value = 205.5
$ws.Range("K55").Value = 205.5
$ws.Range("M55").Value = -32.5
$ws.Range("H61").Value = 2813.25
$ws.Range("I61").Value = 1447.5714
$ws.Range("J61").Value = 5999.8335
$ws.Range("K61").Value = 1447.5714
$ws.Range("L61").Value = 5999.8335
$ws.Range("M61").Value = -1245.5714
$ws.Range("N61").Value = -6403.8335
$ws.Range("H80").Value = 60128
$ws.Range("J80").Value = 60128
$ws.Range("L80").Value = 60128
$ws.Range("N80").Value = -62374
$ws.Range("H83").Value = 60128
$ws.Range("J83").Value = 60128
$ws.Range("L83").Value = 180384
$ws.Range("N83").Value = -191616
$ws.Range("H113").Value = 2813.25
$ws.Range("I113").Value = 1447.5714
$ws.Range("J113").Value = 5999.8335
$ws.Range("K113").Value = 1447.5714
$ws.Range("L113").Value = 5999.8335
$ws.Range("M113").Value = 722.4286
$ws.Range("N113").Value = -10339.8335
$ws.Range("H122").Value = 3706.8538
$ws.Range("I122").Value = 3542.7917
$ws.Range("J122").Value = 3938.4707
$ws.Range("K122").Value = 10628.3751
$ws.Range("L122").Value = 11815.4121
$ws.Range("M122").Value = -8178.375100000001
$ws.Range("N122").Value = -16715.4121
$ws.Range("H126").Value = 2923
$ws.Range("I126").Value = 2718.5715
$ws.Range("K126").Value = 8155.7145
$ws.Range("M126").Value = -5685.7145
$ws.Range("H132").Value = 1443777.2
$ws.Range("I132").Value = 1506289.4
$ws.Range("K132").Value = 4518868.199999999
$ws.Range("M132").Value = -4516338.199999999
$ws.Range("H138").Value = 72000
$ws.Range("J138").Value = 72000
$ws.Range("L138").Value = 72000
$ws.Range("N138").Value = -82280

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 2820.1936
$ws.Range("I113").Value = 1432.7333
$ws.Range("J113").Value = 4120.9375
$ws.Range("K113").Value = 4298.199900000001
$ws.Range("L113").Value = 12362.8125
$ws.Range("M113").Value = -2128.199900000001
$ws.Range("N113").Value = -16702.8125
$ws.Range("H126").Value = 3806.2
$ws.Range("I126").Value = 2826.7273
$ws.Range("J126").Value = 6499.75
$ws.Range("K126").Value = 8480.1819
$ws.Range("L126").Value = 19499.25
$ws.Range("M126").Value = -6010.1819
$ws.Range("N126").Value = -24439.25
$ws.Range("H132").Value = 5298512
$ws.Range("I132").Value = 5752193
$ws.Range("K132").Value = 17256579
$ws.Range("M132").Value = -17254049

Write-Output "Applied scheduled market-data refresh to ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets."
